# Add an 11th "week" (Tube/Group/Individual) of egestion tube ID labels.
# Continues the existing pattern found in column A of Sheet1:
#   E-T{tube}-G{group}-I{individual}-{CA|EP}
# Tube cycles 1-5 (2 groups per tube), Group increments every 4 rows,
# Individual increments by 1 each row, and CA/EP alternate.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 401
$startIndividual = 401
$startGroup = 101

for ($i = 0; $i -lt 40; $i++) {
    $row = $startRow + $i
    $individual = $startIndividual + $i
    $group = $startGroup + [math]::Floor($i / 4)
    $tube = 1 + [math]::Floor($i / 8)
    if ($i % 2 -eq 0) {
        $suffix = "CA"
    } else {
        $suffix = "EP"
    }
    $value = "E-T$tube-G$group-I$individual-$suffix"
    $ws.Cells.Item($row, 1).Value = $value
}
